# Add support for mono camera
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: extend the settings list with new rows 34-36 ---
# Row 33 (commandline mode) is the template for style (label col A = s5, value col B = s7).
$ws1.Range("A33:B33").Copy()
$ws1.Range("A34:B36").PasteSpecial(-4122)   # xlPasteFormats

# Row 36's B value should literally be the text "True" (like B33) rather than a boolean,
# so copy the value from B33 (already stored as shared-string text) instead of assigning it directly.
$ws1.Range("B33").Copy()
$ws1.Range("B36").PasteSpecial(-4163)       # xlPasteValues

$excel.CutCopyMode = 0

# Row 35: new "Overlay mode" setting (color / LRGB / RGB / LHSO / HSO), plus a note in column D.
$ws1.Range("A35").Value = "叠加模式"
$ws1.Range("B35").Value = "color"
$ws1.Range("D35").Value = "彩色相机照片使用color，黑白相机配滤镜照片使用其余。"

# Row 36: new "Use histogram equalization" toggle.
$ws1.Range("A36").Value = "使用直方图均衡化"

# --- Sheet2: add column F, the list source for the new "叠加模式" dropdown ---
$ws2.Range("F1").Value = "color"
$ws2.Range("F2").Value = "LRGB"
$ws2.Range("F3").Value = "RGB"
$ws2.Range("F4").Value = "LHSO"
$ws2.Range("F5").Value = "HSO"

# --- Data validation ---
# B35 gets a list validation sourced from Sheet2!F1:F5.
$ws1.Range("B35").Validation.Add(3, 1, 1, 'Sheet2!$F$1:$F$5')

# B34 and B36 should validate against the same True/False list as B32:B33 (color/mono-style list on Sheet2!B1:B2).
$ws1.Range("B34").Validation.Add(3, 1, 1, 'Sheet2!$B$1:$B$2')
$ws1.Range("B36").Validation.Add(3, 1, 1, 'Sheet2!$B$1:$B$2')

# Sheet2 also picked up an explicit page setup (portrait, paper size 9) in this revision.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- View bookkeeping to mirror the final cursor/viewport position ---
# (Sheet2's own selection moved to its new F2 cell; Sheet1 stays the active/visible tab.)
$ws2.Activate()
$ws2.Range("F2").Select()
$ws1.Activate()
$ws1.Range("A7").Select()
$ws1.Range("D36").Select()
